# This script reshuffles the "weekly" price records (rows 2-21) of the sheet.
# Each data row's own identifying columns (A,B,C,E-K,Q,R,T) stay fixed per
# row, but the variable observation columns (D Fecha, L Calidad, M Volumen,
# N Precio minimo, O Precio maximo, P Precio promedio ponderado, S Precio
# $/Kg) are redistributed among the rows according to a fixed permutation,
# as described by the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of destination row -> source row (values captured BEFORE any
# writes, then assigned to their new destination row).
$mapping = @{
    2  = 20
    3  = 21
    4  = 16
    5  = 8
    6  = 9
    7  = 18
    8  = 19
    9  = 12
    10 = 13
    11 = 11
    12 = 14
    13 = 15
    14 = 2
    15 = 3
    16 = 10
    17 = 4
    18 = 17
    19 = 7
    20 = 5
    21 = 6
}

# Capture the original values for the variable columns of every source row
# before any cell is overwritten.
$original = @{}
foreach ($row in 2..21) {
    $original[$row] = @{
        D = $ws.Cells.Item($row, 4).Value2   # Fecha
        L = $ws.Cells.Item($row, 12).Value2  # Calidad
        M = $ws.Cells.Item($row, 13).Value2  # Volumen
        N = $ws.Cells.Item($row, 14).Value2  # Precio minimo
        O = $ws.Cells.Item($row, 15).Value2  # Precio maximo
        P = $ws.Cells.Item($row, 16).Value2  # Precio promedio ponderado
        S = $ws.Cells.Item($row, 19).Value2  # Precio $/Kg
    }
}

# Apply the captured values to their new destination rows.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $src = $original[$srcRow]

    $ws.Cells.Item($destRow, 4).Value  = $src.D
    $ws.Cells.Item($destRow, 12).Value = $src.L
    $ws.Cells.Item($destRow, 13).Value = $src.M
    $ws.Cells.Item($destRow, 14).Value = $src.N
    $ws.Cells.Item($destRow, 15).Value = $src.O
    $ws.Cells.Item($destRow, 16).Value = $src.P
    $ws.Cells.Item($destRow, 19).Value = $src.S
}
